$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J4 was holding the numeric UI-size value (0.9). Re-fit the UI: shift that
# value down the row as text (K4), bump the old K4 prefab path into L4
# (displacing the stale "MagicBall" placeholder), and zero out J4.

# Preserve L4's prefab-path text by copying K4's current text into it.
$ws.Range("L4").Value = $ws.Range("K4").Value2

# Write "0.899999976158142" into K4 as TEXT (not a number) while keeping
# K4's existing style/number format untouched. A direct .Value assignment
# coerces a numeric-looking string to a real number, so build it as a
# text-literal formula in a scratch cell (no NumberFormat change needed),
# then paste-special just the resulting value across and tidy up.
$scratch = $ws.Range("N4")
$scratch.Formula = '="0.899999976158142"'
$scratch.Copy()
$ws.Range("K4").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("J4").Value = 0
